# método para verificar próxima linha a ser preenchida
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$codigo = 6364025019
$conteudoEmbalagem = "Conteúdo da Embalagem:  1 Luva de lã sintética.  Resistente e prática, não risca e nem danifica a lataria do veículo. Contém punho de algodão para melhor fixação na mão do operador.  Indicada para polimento de superfícies, proporcionando ótimo acabamento e brilho. Ideal para uso automotivo. Atenção: o produto não risca a superfície, mas pode reter resíduos que riscam. Certifique-se de que a luva esteja limpa e livre de resíduos antes de usá-la.  Garantia legal: 90 dias"
$linkFornecedor = "https://www.vonder.com.br/produto/luva_de_l_sinttica_para_polimento_vonder/4468"
$descricaoTitulo = "Luva de lã sintética para polimento VONDER 63.64.025.019"
$conteudoHtml = "<div class=""descricaoProd"">`n               <b>Conteúdo da Embalagem:</b> <br>                     `n                      <p style=""margin-left: 5px; padding-bottom: 10px;"">1 Luva de lã sintética.</p><br>`n                      <p style=""margin-left: 5px; padding-bottom: 10px;"">Resistente e prática, não risca e nem danifica a lataria do veículo. Contém punho de algodão para melhor fixação na mão do operador.</p><br>`n                      <p style=""margin-left: 5px; padding-bottom: 10px;"">Indicada para polimento de superfícies, proporcionando ótimo acabamento e brilho. Ideal para uso automotivo. Atenção: o produto não risca a superfície, mas pode reter resíduos que riscam. Certifique-se de que a luva esteja limpa e livre de resíduos antes de usá-la.</p><br>`n                      <p style=""margin-left: 5px; padding-bottom: 10px;"">Garantia legal: 90 dias</p><br>`n                            <div class=""arquivoItens""> `n                                       </div>`n            </div>"
$detalhesTecnicos = "Comprimento: 250 mm`nLargura: 190 mm`nMaterial: Lã sintética (100% Poliester)`nPunho: Algodão`nMassa aproximada (peso): 0,1 kg"
$categoriaProduto = " Equipamentos de proteção individual, coletiva, sinalização e segurança| Equipamentos para proteção de braços e mãos"

# Linha 2 (dados originais)
$ws.Range("B2").Value = $conteudoEmbalagem
$ws.Range("C2").Value = $linkFornecedor
$ws.Range("D2").Value = $descricaoTitulo
$ws.Range("E2").Value = $conteudoEmbalagem
$ws.Range("F2").Value = $conteudoHtml
$ws.Range("G2").Value = $detalhesTecnicos
$ws.Range("I2").Value = $conteudoHtml
$ws.Range("J2").Value = $categoriaProduto

# Linha 3 (próxima linha a ser preenchida, duplicando os dados)
$ws.Range("A3").Value = $codigo
$ws.Range("B3").Value = $conteudoEmbalagem
$ws.Range("C3").Value = $linkFornecedor
$ws.Range("D3").Value = $descricaoTitulo
$ws.Range("E3").Value = $conteudoEmbalagem
$ws.Range("F3").Value = $conteudoHtml
$ws.Range("G3").Value = $detalhesTecnicos
$ws.Range("I3").Value = $conteudoHtml
$ws.Range("J3").Value = $categoriaProduto

$ws.Range("A1").Select()
